$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column L width from 42 to 50 (stored sheet XML "width" = ColumnWidth + 5/6,
# so request ColumnWidth = 50 - 5/6 to land exactly on a stored width of 50)
$ws.Columns.Item(12).ColumnWidth = 49.166666666666664

# Rename header values in row 1 from input_rowSelection_* to input_rowSelectionCheckbox_*
$ws.Range("A1").Value = "input_rowSelectionCheckbox_class"
$ws.Range("B1").Value = "input_rowSelectionCheckbox_class_1"
$ws.Range("C1").Value = "input_rowSelectionCheckbox_class_2"
$ws.Range("D1").Value = "input_rowSelectionCheckbox_class_3"
$ws.Range("E1").Value = "input_rowSelectionCheckbox_class_4"
$ws.Range("F1").Value = "input_rowSelectionCheckbox_class_5"
$ws.Range("G1").Value = "input_rowSelectionCheckbox_internalRoleRowName"
$ws.Range("H1").Value = "input_rowSelectionCheckbox_internalRoleRowName_1"
$ws.Range("I1").Value = "input_rowSelectionCheckbox_internalRoleRowName_2"
$ws.Range("J1").Value = "input_rowSelectionCheckbox_internalRoleRowName_3"
$ws.Range("K1").Value = "input_rowSelectionCheckbox_internalRoleRowName_4"
$ws.Range("L1").Value = "input_rowSelectionCheckbox_internalRoleRowName_5"
